$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 76666890
$ws.Range("J6").Value = 333333340
$ws.Range("L6").Value = 1000000020
$ws.Range("N6").Value = -1000000244
$ws.Range("H17").Value = 1326201.2
$ws.Range("J17").Value = 1326201.2
$ws.Range("L17").Value = 3978603.6
$ws.Range("N17").Value = -3978939.6
$ws.Range("H19").Value = 530.4737
$ws.Range("I19").Value = 435.63635
$ws.Range("J19").Value = 660.875
$ws.Range("K19").Value = 435.63635
$ws.Range("L19").Value = 660.875
$ws.Range("M19").Value = -260.63635
$ws.Range("N19").Value = -1010.875
$ws.Range("H43").Value = 1104.8182
$ws.Range("J43").Value = 1130.625
$ws.Range("L43").Value = 1130.625
$ws.Range("N43").Value = -1268.625
$ws.Range("H96").Value = 2479
$ws.Range("I96").Value = 1269
$ws.Range("K96").Value = 3807
$ws.Range("M96").Value = -2434
$ws.Range("H106").Value = 3057.3125
$ws.Range("I106").Value = 3005.4614
$ws.Range("K106").Value = 3005.4614
$ws.Range("M106").Value = -2374.4614
$ws.Range("H113").Value = 61229.75
$ws.Range("I113").Value = 73334.31
$ws.Range("K113").Value = 73334.31
$ws.Range("M113").Value = -70080.31
$ws.Range("H132").Value = 8180215
$ws.Range("I132").Value = 9553608
$ws.Range("K132").Value = 28660824
$ws.Range("M132").Value = -28658294
$ws.Range("H133").Value = 200000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 697059.75
$ws.Range("I137").Value = 22489.88
$ws.Range("J137").Value = 2383484.5
$ws.Range("K137").Value = 67469.64
$ws.Range("L137").Value = 7150453.5
$ws.Range("M137").Value = -64919.64
$ws.Range("N137").Value = -7155553.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1475.5294
$ws.Range("J2").Value = 1325.4286
$ws.Range("L2").Value = 1325.4286
$ws.Range("N2").Value = -1551.4286
$ws.Range("H32").Value = 24716.459
$ws.Range("I32").Value = 14315.857
$ws.Range("K32").Value = 14315.857
$ws.Range("M32").Value = -14028.857
$ws.Range("H97").Value = 1283
$ws.Range("J97").Value = 859.5
$ws.Range("L97").Value = 859.5
$ws.Range("N97").Value = -1851.5
$ws.Range("H102").Value = 5875.857
$ws.Range("I102").Value = 7226.2
$ws.Range("K102").Value = 7226.2
$ws.Range("M102").Value = -5604.2
$ws.Range("H116").Value = 1475.5294
$ws.Range("J116").Value = 1325.4286
$ws.Range("L116").Value = 1325.4286
$ws.Range("N116").Value = -5913.4286
$ws.Range("H132").Value = 2536
$ws.Range("I132").Value = 2113.0625
$ws.Range("K132").Value = 6339.1875
$ws.Range("M132").Value = -3809.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1475.5294
$ws.Range("J3").Value = 1325.4286
$ws.Range("L3").Value = 1325.4286
$ws.Range("N3").Value = -1553.4286
$ws.Range("H94").Value = 748.64703
$ws.Range("I94").Value = 700.4167
$ws.Range("K94").Value = 700.4167
$ws.Range("M94").Value = -249.4167
$ws.Range("H134").Value = 1541.1111
$ws.Range("I134").Value = 1134.6129
$ws.Range("J134").Value = 4061.4
$ws.Range("K134").Value = 3403.8387
$ws.Range("L134").Value = 12184.2
$ws.Range("M134").Value = -868.8387000000002
$ws.Range("N134").Value = -17254.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8071825
$ws.Range("I31").Value = 3015946
$ws.Range("K31").Value = 3015946
$ws.Range("M31").Value = -3015651
$ws.Range("H34").Value = 8071825
$ws.Range("I34").Value = 3015946
$ws.Range("K34").Value = 3015946
$ws.Range("M34").Value = -3015744
$ws.Range("H58").Value = 2537.9688
$ws.Range("I58").Value = 2048.4783
$ws.Range("K58").Value = 2048.4783
$ws.Range("M58").Value = -1845.4783
$ws.Range("H136").Value = 2537.9688
$ws.Range("I136").Value = 2048.4783
$ws.Range("K136").Value = 6145.4349
$ws.Range("M136").Value = -3595.4349
$ws.Range("H140").Value = 75489.25
$ws.Range("J140").Value = 75489.25
$ws.Range("L140").Value = 75489.25
$ws.Range("N140").Value = -85849.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 3402.6
$ws.Range("I59").Value = 1903.25
$ws.Range("J59").Value = 9400
$ws.Range("K59").Value = 5709.75
$ws.Range("L59").Value = 28200
$ws.Range("M59").Value = -5169.75
$ws.Range("N59").Value = -29280
$ws.Range("H80").Value = 1500
$ws.Range("H83").Value = 1500
$ws.Range("H92").Value = 2000
$ws.Range("J92").Value = 2000
$ws.Range("L92").Value = 6000
$ws.Range("N92").Value = -8496
$ws.Range("H107").Value = 27778596
$ws.Range("I107").Value = 37037496
$ws.Range("J107").Value = 1894.6666
$ws.Range("K107").Value = 111112488
$ws.Range("L107").Value = 5683.9998
$ws.Range("M107").Value = -111110568
$ws.Range("N107").Value = -9523.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H80").Value = 250001060
$ws.Range("J80").Value = 1400
$ws.Range("L80").Value = 1400
$ws.Range("N80").Value = -3396
$ws.Range("H83").Value = 250001060
$ws.Range("J83").Value = 1400
$ws.Range("L83").Value = 7000
$ws.Range("N83").Value = -16984
$ws.Range("H97").Value = 1385.375
$ws.Range("I97").Value = 1426.0938
$ws.Range("K97").Value = 1426.0938
$ws.Range("M97").Value = -930.0938000000001
$ws.Range("H107").Value = 1089.3334
$ws.Range("I107").Value = 941.2
$ws.Range("K107").Value = 941.2
$ws.Range("M107").Value = 978.8
$ws.Range("H113").Value = 4561.769
$ws.Range("I113").Value = 4163
$ws.Range("K113").Value = 4163
$ws.Range("M113").Value = -1993
$ws.Range("H122").Value = 7873.735
$ws.Range("J122").Value = 3276.923
$ws.Range("L122").Value = 9830.769
$ws.Range("N122").Value = -14730.769
$ws.Range("H123").Value = 63949.8
$ws.Range("J123").Value = 73449.75
$ws.Range("L123").Value = 73449.75
$ws.Range("N123").Value = -78349.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4120.409
$ws.Range("I40").Value = 3070.8333
$ws.Range("K40").Value = 3070.8333
$ws.Range("M40").Value = -2934.8333
$ws.Range("H61").Value = 4081.6428
$ws.Range("I61").Value = 3137.2222
$ws.Range("K61").Value = 3137.2222
$ws.Range("M61").Value = -2935.2222
$ws.Range("H93").Value = 1590738.8
$ws.Range("I93").Value = 1855195.1
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1855195.1
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -1853947.1
$ws.Range("N93").Value = -6496
$ws.Range("H113").Value = 4081.6428
$ws.Range("I113").Value = 3137.2222
$ws.Range("K113").Value = 3137.2222
$ws.Range("M113").Value = -967.2222000000002
$ws.Range("H132").Value = 4583.773
$ws.Range("I132").Value = 3810.3076
$ws.Range("J132").Value = 5701
$ws.Range("K132").Value = 11430.9228
$ws.Range("L132").Value = 17103
$ws.Range("M132").Value = -8900.9228
$ws.Range("N132").Value = -22163
$ws.Range("H137").Value = 67475
$ws.Range("I137").Value = 30000
$ws.Range("J137").Value = 79966.664
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 79966.664
$ws.Range("M137").Value = -24900
$ws.Range("N137").Value = -90166.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8213.736999999999
$ws.Range("I81").Value = 15148
$ws.Range("J81").Value = 4168.75
$ws.Range("K81").Value = 30296
$ws.Range("L81").Value = 8337.5
$ws.Range("M81").Value = -29235
$ws.Range("N81").Value = -10459.5
$ws.Range("H84").Value = 8213.736999999999
$ws.Range("I84").Value = 15148
$ws.Range("J84").Value = 4168.75
$ws.Range("K84").Value = 151480
$ws.Range("L84").Value = 41687.5
$ws.Range("M84").Value = -146176
$ws.Range("N84").Value = -52295.5
$ws.Range("H122").Value = 12099
$ws.Range("I122").Value = 4156.4287
$ws.Range("J122").Value = 21365.334
$ws.Range("K122").Value = 12469.2861
$ws.Range("L122").Value = 64096.00199999999
$ws.Range("M122").Value = -10019.2861
$ws.Range("N122").Value = -68996.00199999999
$ws.Range("H126").Value = 2904
$ws.Range("I126").Value = 2554.389
$ws.Range("K126").Value = 7663.167
$ws.Range("M126").Value = -5193.167
$ws.Range("H132").Value = 9695.362999999999
$ws.Range("I132").Value = 9695.362999999999
$ws.Range("K132").Value = 29086.089
$ws.Range("M132").Value = -26556.089
$ws.Range("H139").Value = 106716.664
$ws.Range("J139").Value = 106716.664
$ws.Range("L139").Value = 106716.664
$ws.Range("N139").Value = -116996.664
